$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.563.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3598'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.86%  '
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9000'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07782'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.834.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.286'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.336'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008562'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.603.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.979'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.042.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.16%  '
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.060'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.870'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08716'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.126'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7531'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.753'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.436'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05112'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.906'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("E42").Value = '  +2.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.767'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1508'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.062'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4738'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.01%  '
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.998'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.581'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05987'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
